$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(44992, 1040, 2000, 2500, 2250, 750)
    3 = @(45084, 900, 2000, 2500, 2250, 750)
    4 = @(45006, 1100, 2000, 2500, 2250, 750)
    5 = @(45077, 760, 2000, 2500, 2250, 750)
    6 = @(45035, 1100, 2000, 2500, 2250, 750)
    7 = @(44685, 400, 1500, 2000, 1750, 583)
    8 = @(45034, 1100, 2000, 2500, 2250, 750)
    9 = @(44911, 700, 1800, 2000, 1900, 633)
    10 = @(44951, 800, 2000, 2500, 2250, 750)
    11 = @(44965, 1120, 2000, 2500, 2250, 750)
    12 = @(45028, 1000, 2000, 2500, 2250, 750)
    13 = @(44827, 1200, 2000, 2500, 2250, 750)
    14 = @(44883, 500, 1800, 2000, 1900, 633)
    15 = @(44970, 800, 2000, 2500, 2250, 750)
    16 = @(44848, 1000, 1500, 2000, 1750, 583)
    17 = @(45007, 1160, 2000, 2500, 2250, 750)
    18 = @(44971, 1000, 2000, 2500, 2250, 750)
    19 = @(44985, 1000, 2000, 2500, 2250, 750)
    20 = @(44964, 1000, 2000, 2500, 2250, 750)
    21 = @(45020, 1200, 2000, 2500, 2250, 750)
    22 = @(44978, 1000, 1800, 2000, 1900, 633)
    23 = @(45091, 800, 2000, 2500, 2250, 750)
    24 = @(44953, 1000, 2000, 2500, 2250, 750)
    25 = @(45062, 1100, 2000, 2500, 2250, 750)
    26 = @(44910, 1000, 1800, 2000, 1900, 633)
    27 = @(45041, 1160, 2000, 2500, 2250, 750)
    28 = @(45070, 800, 2000, 2500, 2250, 750)
    29 = @(44881, 500, 1900, 2000, 1950, 650)
    30 = @(44999, 1100, 2000, 2500, 2250, 750)
    31 = @(45013, 1100, 2000, 2500, 2250, 750)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]  # P - Precio $/Kg
}
